$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text in Q1
$ws.Range("Q1").Value = "wtkappa.scale_trim"

# Flip signs on E2 and F2
$ws.Range("E2").Value = 0.02351246133036377
$ws.Range("F2").Value = -0.01081937260331701

# Update Q2 value
$ws.Range("Q2").Value = 0.7808705382933501
